# Insert a new "Industry" column between "Stock Name" (B) and "Mutual Fund" (C),
# shifting all subsequent columns one position to the right, then populate the
# new column with each holding's industry classification.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at position C; existing C..I shift to D..J.
$ws.Columns.Item(3).Insert()

# Header for the new column.
$ws.Range("C1").Value = "Industry"

# Industry values for rows 2-23 (aligned with ISIN / Stock Name in A/B).
$industries = @(
    "Finance",                   # 2  Piramal Finance Ltd
    "Finance",                   # 3  Shriram Finance Limited
    "Finance",                   # 4  Capri Global Capital Limited
    "Banks",                     # 5  HDFC Bank Limited
    "Finance",                   # 6  LIC Housing Finance Ltd
    "Banks",                     # 7  Kotak Mahindra Bank Limited
    "Finance",                   # 8  Bajaj Finance Limited
    "Capital Markets",           # 9  ICICI Prudential AMC Ltd
    "Capital Markets",           # 10 HDFC Asset Management Company Ltd
    "Insurance",                 # 11 HDFC Life Insurance Co Ltd
    "Metals & Minerals Trading", # 12 Adani Enterprises Limited
    "Capital Markets",           # 13 Nippon Life India Asset Management Ltd
    "Banks",                     # 14 ICICI Bank Limited
    "Insurance",                 # 15 ICICI Prudential Life Insurance Co Ltd
    "Finance",                   # 16 SBI Cards & Payment Services Ltd
    "Capital Markets",           # 17 Anand Rathi Wealth Limited
    "Banks",                     # 18 Kotak Mahindra Bank Limited
    "Insurance",                 # 19 Canara HSBC Life Insurance Company Ltd
    "Insurance",                 # 20 SBI Life Insurance Company Limited
    "Insurance",                 # 21 Life Insurance Corporation Of India
    "Banks",                     # 22 State Bank of India
    "Finance"                    # 23 Bajaj Finserv Ltd.
)

for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $industries[$i]
}
